# #5: cash & deposit done
# Add bank/deposit metadata columns (G:M) to the 存款 (deposit) sheet,
# mirroring the property_category/category/date/legislator_name/
# legislator_id/source_file/index columns already present on the other
# asset sheets (land/building/car/...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1): new columns G-M ---
$ws.Cells.Item(1, 7).Value  = "property_category"
$ws.Cells.Item(1, 8).Value  = "category"
$ws.Cells.Item(1, 9).Value  = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

$headerRange = $ws.Range("G1:M1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows 2-12: new columns G-M ---
# row number -> original "index" column value (matches column A on this sheet)
$rows = @(
    @(2, 47),
    @(3, 48),
    @(4, 49),
    @(5, 50),
    @(6, 51),
    @(7, 52),
    @(8, 53),
    @(9, 54),
    @(10, 55),
    @(11, 56),
    @(12, 57)
)

foreach ($item in $rows) {
    $r = $item[0]
    $idx = $item[1]

    $ws.Cells.Item($r, 7).Value  = "deposit"
    $ws.Cells.Item($r, 8).Value  = "normal"
    $ws.Cells.Item($r, 9).Value  = "2013-12-02"
    $ws.Cells.Item($r, 10).Value = "李俊俋"
    $ws.Cells.Item($r, 11).Value = 1738
    $ws.Cells.Item($r, 12).Value = "tmp52b51"
    $ws.Cells.Item($r, 13).Value = $idx
}
